# Update "想去人数" (people interested) counts that changed between scrapes.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F14").Value = 128
$ws1.Range("F19").Value = 4102
$ws1.Range("F23").Value = 532
$ws1.Range("F24").Value = 1644
$ws1.Range("F28").Value = 2198
$ws1.Range("F47").Value = 602
$ws1.Range("F48").Value = 713

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 567

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F16").Value = 128
$ws4.Range("F21").Value = 4102
$ws4.Range("F27").Value = 532
$ws4.Range("F28").Value = 1644
$ws4.Range("F32").Value = 2198
$ws4.Range("F49").Value = 602
$ws4.Range("F50").Value = 713
